$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update all C2:C402 values from 0.1 to 0.05
$ws.Range("C2:C402").Value = 0.05

# Reset the selection to C1 (single cell) instead of C2:C402
$ws.Range("C1").Select()
